$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 1491.8795
$ws.Range("I15").Value = 1491.8795
$ws.Range("K15").Value = 4475.6385
$ws.Range("M15").Value = -4306.6385

$ws.Range("H40").Value = 1228.88
$ws.Range("I40").Value = 1125.125
$ws.Range("J40").Value = 1413.3334
$ws.Range("K40").Value = 1125.125
$ws.Range("L40").Value = 1413.3334
$ws.Range("M40").Value = -950.125
$ws.Range("N40").Value = -1763.3334

$ws.Range("H51").Value = 2640.6155
$ws.Range("I51").Value = 2091.5
$ws.Range("J51").Value = 2884.6667
$ws.Range("K51").Value = 2091.5
$ws.Range("L51").Value = 2884.6667
$ws.Range("M51").Value = -1607.5
$ws.Range("N51").Value = -3852.6667

$ws.Range("H63").Value = 10000
$ws.Range("J63").Value = 10000
$ws.Range("L63").Value = 10000
$ws.Range("N63").Value = -11248

$ws.Range("H66").Value = 10000
$ws.Range("J66").Value = 10000
$ws.Range("L66").Value = 30000
$ws.Range("N66").Value = -36240

$ws.Range("H100").Value = 12822147
$ws.Range("I100").Value = 17544776
$ws.Range("J100").Value = 3584.2856
$ws.Range("K100").Value = 17544776
$ws.Range("L100").Value = 3584.2856
$ws.Range("M100").Value = -17544235
$ws.Range("N100").Value = -4666.2856

$ws.Range("H137").Value = 886.6591
$ws.Range("I137").Value = 782.075
$ws.Range("J137").Value = 1932.5
$ws.Range("K137").Value = 2346.225
$ws.Range("L137").Value = 5797.5
$ws.Range("M137").Value = 203.7749999999996
$ws.Range("N137").Value = -10897.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1871
$ws.Range("I2").Value = 763.3077
$ws.Range("J2").Value = 4271
$ws.Range("K2").Value = 763.3077
$ws.Range("L2").Value = 4271
$ws.Range("M2").Value = -650.3077
$ws.Range("N2").Value = -4497

$ws.Range("H5").Value = 1500.2
$ws.Range("I5").Value = 1500.2
$ws.Range("K5").Value = 1500.2
$ws.Range("M5").Value = -1388.2

$ws.Range("H32").Value = 1576.84
$ws.Range("I32").Value = 1403.5333
$ws.Range("J32").Value = 3136.6
$ws.Range("K32").Value = 1403.5333
$ws.Range("L32").Value = 3136.6
$ws.Range("M32").Value = -1116.5333
$ws.Range("N32").Value = -3710.6

$ws.Range("H116").Value = 1871
$ws.Range("I116").Value = 763.3077
$ws.Range("J116").Value = 4271
$ws.Range("K116").Value = 763.3077
$ws.Range("L116").Value = 4271
$ws.Range("M116").Value = 1530.6923
$ws.Range("N116").Value = -8859

$ws.Range("H122").Value = 5684042
$ws.Range("I122").Value = 2232.6155
$ws.Range("J122").Value = 13891100
$ws.Range("K122").Value = 6697.8465
$ws.Range("L122").Value = 41673300
$ws.Range("M122").Value = -4247.8465
$ws.Range("N122").Value = -41678200

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1871
$ws.Range("I3").Value = 763.3077
$ws.Range("J3").Value = 4271
$ws.Range("K3").Value = 763.3077
$ws.Range("L3").Value = 4271
$ws.Range("M3").Value = -649.3077
$ws.Range("N3").Value = -4499

$ws.Range("H4").Value = 1500.2
$ws.Range("I4").Value = 1500.2
$ws.Range("K4").Value = 1500.2
$ws.Range("M4").Value = -1385.2

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1100
$ws.Range("I16").Value = 850
$ws.Range("J16").Value = 1350
$ws.Range("K16").Value = 850
$ws.Range("L16").Value = 1350
$ws.Range("M16").Value = -563
$ws.Range("N16").Value = -1924

$ws.Range("H113").Value = 1100
$ws.Range("I113").Value = 850
$ws.Range("J113").Value = 1350
$ws.Range("K113").Value = 850
$ws.Range("L113").Value = 1350
$ws.Range("M113").Value = 1320
$ws.Range("N113").Value = -5690

$ws.Range("H122").Value = 7693100.5
$ws.Range("I122").Value = 753.6
$ws.Range("J122").Value = 18182664
$ws.Range("K122").Value = 2260.8
$ws.Range("L122").Value = 54547992
$ws.Range("M122").Value = 189.1999999999998
$ws.Range("N122").Value = -54552892

$ws.Range("H132").Value = 4880534
$ws.Range("I132").Value = 2109.9033
$ws.Range("J132").Value = 20003648
$ws.Range("K132").Value = 6329.7099
$ws.Range("L132").Value = 60010944
$ws.Range("M132").Value = -3799.7099
$ws.Range("N132").Value = -60016004

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 1016.5143
$ws.Range("I68").Value = 770.7143
$ws.Range("J68").Value = 1180.381
$ws.Range("K68").Value = 2312.1429
$ws.Range("L68").Value = 3541.143
$ws.Range("M68").Value = -1501.1429
$ws.Range("N68").Value = -5163.143

$ws.Range("H71").Value = 1016.5143
$ws.Range("I71").Value = 770.7143
$ws.Range("J71").Value = 1180.381
$ws.Range("K71").Value = 6936.428699999999
$ws.Range("L71").Value = 10623.429
$ws.Range("M71").Value = -2880.428699999999
$ws.Range("N71").Value = -18735.429

$ws.Range("H117").Value = 112287.22
$ws.Range("I117").Value = 917.4
$ws.Range("J117").Value = 251499.5
$ws.Range("K117").Value = 2752.2
$ws.Range("L117").Value = 754498.5
$ws.Range("M117").Value = 689.8000000000002
$ws.Range("N117").Value = -761382.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H63").Value = 29500
$ws.Range("J63").Value = 29500
$ws.Range("L63").Value = 29500
$ws.Range("N63").Value = -30872

$ws.Range("H66").Value = 29500
$ws.Range("J66").Value = 29500
$ws.Range("L66").Value = 88500
$ws.Range("N66").Value = -95364

$ws.Range("H113").Value = 2318.0908
$ws.Range("I113").Value = 2249.8333
$ws.Range("J113").Value = 2400
$ws.Range("K113").Value = 2249.8333
$ws.Range("L113").Value = 2400
$ws.Range("M113").Value = -79.83329999999978
$ws.Range("N113").Value = -6740

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H41").Value = 4082.5
$ws.Range("I41").Value = 1033
$ws.Range("K41").Value = 1033
$ws.Range("M41").Value = -595

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H64").Value = 0
$ws.Range("J64").Value = 0
$ws.Range("L64").Value = 0
$ws.Range("N64").Value = $null

$ws.Range("H67").Value = 0
$ws.Range("J67").Value = 0
$ws.Range("L67").Value = 0
$ws.Range("N67").Value = $null

$ws.Range("H132").Value = 72586930
$ws.Range("I132").Value = 160716930
$ws.Range("J132").Value = 9281.058999999999
$ws.Range("K132").Value = 482150790
$ws.Range("L132").Value = 27843.177
$ws.Range("M132").Value = -482148260
$ws.Range("N132").Value = -32903.177
